# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for column G, rows 2-29
$kValues = @{
    2  = 7
    3  = 6
    4  = 7
    5  = 5
    6  = 4
    7  = 5
    8  = 1
    9  = 0
    10 = 5
    11 = 4
    12 = 6
    13 = 6
    14 = 3
    15 = 5
    16 = 6
    17 = 3
    18 = 7
    19 = 2
    20 = 8
    21 = 6
    22 = 5
    23 = 3
    24 = 7
    25 = 4
    26 = 3
    27 = 1
    28 = 3
    29 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
